$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "(01-2023 to 01-2024)" -> "(08-2022 to 01-2024)" in the Associate
#    DevOps Engineer / CureMD entry. Word originally stored this as two
#    runs: "01" | "-2023 to 01-2024)". The edit only changes two single
#    characters ("1"->"8" and "3"->"2"), but Word's interactive editing
#    re-split the surrounding text into five runs. We locate the exact
#    range with Find, then use temporary bookmarks to pin the run-split
#    points (bookmarks stop the engine's run auto-coalescing) before
#    doing the two single-character replacements, then remove the
#    temporary bookmarks again.
# ---------------------------------------------------------------------
$dateRange = $d.Content
$dateFound = $dateRange.Find.Execute("01-2023 to 01-2024)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $dateFound) {
    throw "Could not locate the '01-2023 to 01-2024)' date text"
}
$dateStart = $dateRange.Start

$d.Bookmarks.Add("zzSplit1", $d.Range($dateStart + 1, $dateStart + 1)) | Out-Null
$d.Bookmarks.Add("zzSplit2", $d.Range($dateStart + 2, $dateStart + 2)) | Out-Null
$d.Bookmarks.Add("zzSplit3", $d.Range($dateStart + 6, $dateStart + 6)) | Out-Null
$d.Bookmarks.Add("zzSplit4", $d.Range($dateStart + 7, $dateStart + 7)) | Out-Null

$d.Range($dateStart + 1, $dateStart + 2).Text = "8"
$d.Range($dateStart + 6, $dateStart + 7).Text = "2"

$d.Bookmarks.Item("zzSplit1").Delete()
$d.Bookmarks.Item("zzSplit2").Delete()
$d.Bookmarks.Item("zzSplit3").Delete()
$d.Bookmarks.Item("zzSplit4").Delete()

# ---------------------------------------------------------------------
# 2) Split "Implemented S3 bucket policies and encryption for secure
#    data storage." into "Implemented S3 buck" + "et policies..." with
#    the document's "_GoBack" bookmark (last-edit-location marker)
#    sitting at the split point. This also moves "_GoBack" away from
#    its old location at the end of the "Version Control Systems:" run.
# ---------------------------------------------------------------------
try {
    $d.Bookmarks.Item("_GoBack").Delete()
} catch {
}

$s3Range = $d.Content
$s3Found = $s3Range.Find.Execute("Implemented S3 bucket policies and encryption for secure data storage.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $s3Found) {
    throw "Could not locate the S3 bucket bullet text"
}
$s3SplitPos = $s3Range.Start + "Implemented S3 buck".Length

$d.Bookmarks.Add("_GoBack", $d.Range($s3SplitPos, $s3SplitPos)) | Out-Null

# Rewrite the trailing half through a throwaway placeholder so the
# engine mints a brand-new run (without the original run's rsid
# attributes) for the second half, matching how Word splits a run when
# text is actually (re)typed at the cursor.
$tailRange = $d.Range($s3SplitPos, $s3SplitPos + "et policies and encryption for secure data storage.".Length)
$tailRange.Text = "@@PLACEHOLDER@@"
$tailRange2 = $d.Range($s3SplitPos, $s3SplitPos + "@@PLACEHOLDER@@".Length)
$tailRange2.Text = "et policies and encryption for secure data storage."

Write-Output "done"
